$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-10 from 45243 to 45244
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
